$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.602.40"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "2.612.18"
$ws.Range("E3").Value = "  +4.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.54"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.72"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.532"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.556"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.43"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.67"
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0817"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.24"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "3.022.85"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("D16").Value = "2.631.00"
$ws.Range("E16").Value = "  +4.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.870"
$ws.Range("E17").Value = "  +3.06%  "
$ws.Range("D18").Value = "49.568.27"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.08"
$ws.Range("E19").Value = "  +10.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.29"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "0.0₃0948"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.48"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "278.86"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.42"
$ws.Range("E26").Value = "  +2.50%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.94"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.34"
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.80"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.68"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.41"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0790"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.04"
$ws.Range("E37").Value = "  +4.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.72"
$ws.Range("E38").Value = "  +1.59%  "
$ws.Range("E39").Value = "  +5.95%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.112"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.95"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.58"
$ws.Range("E42").Value = "  +4.87%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0314"
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.33"
$ws.Range("E45").Value = "  +5.07%  "
$ws.Range("D46").Value = "2.046.83"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.23"
$ws.Range("E47").Value = "  +12.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("E48").Value = "  +8.48%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.34"
$ws.Range("E50").Value = "  +3.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.63"
$ws.Range("E51").Value = "  +1.15%  "
